$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marks: "Marking" row total right answers (B11) and "Total" row right answers (B12)
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 120

# Update the fraction display in E12 (correct/total marks)
$ws.Range("E12").Value = "120/140"
